$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(44277, "", 6570, "Pris banktjänster enligt faktura", 60, "")
    ,@(44277, "", "", "Pris banktjänster enligt faktura", 0, "")
    ,@(44277, "", 1930, "Pris banktjänster enligt faktura", "", 60)
    ,@(44277, "", 4010, "NGROCERIES K0135", 355.36, "")
    ,@(44277, "", 2645, "NGROCERIES K0135", 42.64, "")
    ,@(44277, "", 1930, "NGROCERIES K0135", "", 398)
    ,@(44277, "", 4010, "SNABBGROSS SOLNA K0135", 278, "")
    ,@(44277, "", 2645, "SNABBGROSS SOLNA K0135", 33.36, "")
    ,@(44277, "", 1930, "SNABBGROSS SOLNA K0135", "", 311.36)
    ,@(44279, "Reko59", 3011, "Reko Swish +46725248271", "", 1062.5)
    ,@(44279, "Reko59", 2611, "Reko Swish +46725248271", "", 127.5)
    ,@(44279, "Reko59", 1930, "Reko Swish +46725248271", 1190, "")
    ,@(44279, "", 4010, "SNABBGROSS SOLNA K0135", 657.54, "")
    ,@(44279, "", 2645, "SNABBGROSS SOLNA K0135", 78.9, "")
    ,@(44279, "", 1930, "SNABBGROSS SOLNA K0135", "", 736.44)
    ,@(44280, "Reko60", 3011, "Reko Swish +46700376635", "", 345.54)
    ,@(44280, "Reko60", 2611, "Reko Swish +46700376635", "", 41.46)
    ,@(44280, "Reko60", 1930, "Reko Swish +46700376635", 387, "")
    ,@(44280, "8251020", 3011, "Order 8251020 Swish +46707678891", "", 588.39)
    ,@(44280, "8251020", 2611, "Order 8251020 Swish +46707678891", "", 70.61)
    ,@(44280, "8251020", 1930, "Order 8251020 Swish +46707678891", 659, "")
    ,@(44280, "", 6400, "VISTAPRINT K0135", 844, "")
    ,@(44280, "", 2641, "VISTAPRINT K0135", 211, "")
    ,@(44280, "", 1930, "VISTAPRINT K0135", "", 1055)
    ,@(44282, "", 4010, "SNABBGROSS SOLNA K0135", 1333.05, "")
    ,@(44282, "", 2645, "SNABBGROSS SOLNA K0135", 159.97, "")
    ,@(44282, "", 1930, "SNABBGROSS SOLNA K0135", "", 1493.02)
)

$startRow = 479
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Column B: receipt / reko / order number references are stored as text in
    # this sheet. Purely-numeric-looking values need to be forced to text so they
    # don't get auto-converted to numbers; leading apostrophe forces text entry,
    # then the style is reset back to Normal so no stray number-format is left on
    # the cell (matching the rest of the sheet, which uses the default style).
    $bVal = $row[1]
    $bCell = $ws.Cells.Item($r, 2)
    if ($bVal -ne "" -and $bVal -match "^[0-9]+$") {
        $bCell.Value = "'" + $bVal
        $bCell.Style = "Normal"
    } else {
        $bCell.Value = $bVal
    }

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
